$wb = $excel.ActiveWorkbook

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1899.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1899.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = ""
$ws.Range("M86").Value = 1899.5
$ws.Range("N86").Value = -4145.5

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1899.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1899.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = ""
$ws.Range("M89").Value = 9497.5
$ws.Range("N89").Value = -20729.5

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 655.9091
$ws.Range("I98").Value = 655.9091
$ws.Range("K98").Value = 655.9091
$ws.Range("M98").Value = 842.0909

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2679.2974
$ws.Range("J112").Value = 3070.8965
$ws.Range("L112").Value = 9212.6895
$ws.Range("N112").Value = -11428.6895

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7531.8887
$ws.Range("I113").Value = 7531.8887
$ws.Range("K113").Value = 7531.8887
$ws.Range("M113").Value = -4277.8887

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 655.9091
$ws.Range("I122").Value = 655.9091
$ws.Range("K122").Value = 1967.7273
$ws.Range("M122").Value = 482.2727

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3327.5454
$ws.Range("I141").Value = 2416.5
$ws.Range("J141").Value = 4420.8
$ws.Range("K141").Value = 7249.5
$ws.Range("L141").Value = 13262.4
$ws.Range("M141").Value = -2069.5
$ws.Range("N141").Value = -23622.4

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17039.197
$ws.Range("I32").Value = 6154.104
$ws.Range("J32").Value = 57230.31
$ws.Range("K32").Value = 6154.104
$ws.Range("L32").Value = 57230.31
$ws.Range("M32").Value = -5867.104
$ws.Range("N32").Value = -57804.31

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 654.6
$ws.Range("I74").Value = 654.6
$ws.Range("K74").Value = 654.6
$ws.Range("M74").Value = 219.4

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 654.6
$ws.Range("I77").Value = 654.6
$ws.Range("K77").Value = 3273
$ws.Range("M77").Value = 1095

# Sheet ARM, row 92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 23775
$ws.Range("J92").Value = 23775
$ws.Range("L92").Value = 23775
$ws.Range("N92").Value = -28767

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1398.0952
$ws.Range("I97").Value = 1408.8889
$ws.Range("K97").Value = 1408.8889
$ws.Range("M97").Value = -912.8888999999999

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2119.389
$ws.Range("I122").Value = 2122.1875
$ws.Range("J122").Value = 2097
$ws.Range("K122").Value = 6366.5625
$ws.Range("L122").Value = 6291
$ws.Range("M122").Value = -3916.5625
$ws.Range("N122").Value = -11191

# Sheet BSM, row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 46015.75
$ws.Range("J35").Value = 46015.75
$ws.Range("L35").Value = 46015.75
$ws.Range("N35").Value = -46635.75

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4377.778
$ws.Range("I86").Value = 4377.778
$ws.Range("K86").Value = 4377.778
$ws.Range("M86").Value = -3254.778

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4377.778
$ws.Range("I89").Value = 4377.778
$ws.Range("K89").Value = 21888.89
$ws.Range("M89").Value = -16272.89

# Sheet CRP, row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 14142
$ws.Range("J88").Value = 14142
$ws.Range("L88").Value = 14142
$ws.Range("N88").Value = -14954

# Sheet CRP, row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 14142
$ws.Range("J91").Value = 14142
$ws.Range("L91").Value = 14142
$ws.Range("N91").Value = -16950

# Sheet CRP, row 92
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 29999.5
$ws.Range("J92").Value = 29999.5
$ws.Range("L92").Value = 29999.5
$ws.Range("N92").Value = -34991.5

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3788.2222
$ws.Range("I122").Value = 3936.875
$ws.Range("J122").Value = 2599
$ws.Range("K122").Value = 11810.625
$ws.Range("L122").Value = 7797
$ws.Range("M122").Value = -9360.625
$ws.Range("N122").Value = -12697

# Sheet CUL, row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 325.5
$ws.Range("J34").Value = 309.4
$ws.Range("L34").Value = 928.1999999999999
$ws.Range("N34").Value = -1096.2

# Sheet CUL, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4750
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 144906.58
$ws.Range("I122").Value = 168832.67
$ws.Range("K122").Value = 506498.01
$ws.Range("M122").Value = -504048.01

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3384.4187
$ws.Range("I40").Value = 2581.5557
$ws.Range("K40").Value = 2581.5557
$ws.Range("M40").Value = -2445.5557

# Sheet LTW, row 56
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = ""

# Sheet LTW, row 64
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 39147.145
$ws.Range("J64").Value = 39147.145
$ws.Range("L64").Value = 39147.145
$ws.Range("N64").Value = -39597.145

# Sheet LTW, row 67
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 39147.145
$ws.Range("J67").Value = 39147.145
$ws.Range("L67").Value = 39147.145
$ws.Range("N67").Value = -40707.145

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2903.3914
$ws.Range("I132").Value = 2915.4285
$ws.Range("J132").Value = 2884.6667
$ws.Range("K132").Value = 8746.2855
$ws.Range("L132").Value = 8654.000100000001
$ws.Range("M132").Value = -6216.2855
$ws.Range("N132").Value = -13714.0001

# Sheet WVR, row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 33983
$ws.Range("J68").Value = 33983
$ws.Range("L68").Value = 33983
$ws.Range("N68").Value = -35605

# Sheet WVR, row 69
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 10369.667
$ws.Range("J69").Value = 10369.667
$ws.Range("L69").Value = 10369.667
$ws.Range("N69").Value = -11867.667

# Sheet WVR, row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 33983
$ws.Range("J71").Value = 33983
$ws.Range("L71").Value = 101949
$ws.Range("N71").Value = -110061

# Sheet WVR, row 72
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 10369.667
$ws.Range("J72").Value = 10369.667
$ws.Range("L72").Value = 31109.001
$ws.Range("N72").Value = -38597.001

# Sheet WVR, row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 26899.5
$ws.Range("J80").Value = 26899.5
$ws.Range("L80").Value = 26899.5
$ws.Range("N80").Value = -28895.5

# Sheet WVR, row 82
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 47283.6
$ws.Range("J82").Value = 47283.6
$ws.Range("L82").Value = 47283.6
$ws.Range("N82").Value = -48049.6

# Sheet WVR, row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 26899.5
$ws.Range("J83").Value = 26899.5
$ws.Range("L83").Value = 80698.5
$ws.Range("N83").Value = -90682.5

# Sheet WVR, row 85
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 47283.6
$ws.Range("J85").Value = 47283.6
$ws.Range("L85").Value = 47283.6
$ws.Range("N85").Value = -49935.6
